$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A87").Value2 = 'let email = document.getElementsByName("email")[0].value;'
$ws.Range("B87").Value2 = $ws.Range("B86").Value2

$ws.Range("A87").Select()
